$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37; this shifts existing rows 37-104 down to 38-105,
# matching the dimension change from A1:T104 to A1:T105.
$ws.Rows.Item(37).EntireRow.Insert()

# Populate the newly inserted row 37 with the new weekly record.
$ws.Range("A37").Value = 10
$ws.Range("B37").Value = "Vega Modelo de Temuco"
$ws.Range("C37").Value = "La Araucanía"
$ws.Range("D37").Value = (Get-Date -Year 2023 -Month 8 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E37").Value = 9
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100107
$ws.Range("H37").Value = "Otros"
$ws.Range("I37").Value = 100107011
$ws.Range("J37").Value = "Tuna"
$ws.Range("K37").Value = "Sin especificar"
$ws.Range("L37").Value = "Especial"
$ws.Range("M37").Value = 45
$ws.Range("N37").Value = 32000
$ws.Range("O37").Value = 32000
$ws.Range("P37").Value = 32000
$ws.Range("Q37").Value = "$/caja 16 kilos"
$ws.Range("R37").Value = "Provincia de Los Andes"
$ws.Range("S37").Value = 2000
$ws.Range("T37").Value = 16
